$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Natuurkunde", "NA", 1),
    @("Scheikunde", "SK", 1),
    @("Frans", "FA", 1),
    @("Informatica", "IN", 0)
)

$r = 13
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("A17").Select()
